$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.672.91"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "3.335.14"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'581.94"
$ws.Range("E5").Value = "  -0.86%  "
$ws.Range("D6").Value = "'176.21"
$ws.Range("E6").Value = "  -3.61%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'0.589"
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "3.333.05"
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("D10").Value = "'0.178"
$ws.Range("E10").Value = "  +0.13%  "
$ws.Range("D11").Value = "'0.577"
$ws.Range("E11").Value = "  -0.50%  "
$ws.Range("D12").Value = "'45.45"
$ws.Range("E12").Value = "  -1.91%  "
$ws.Range("E13").Value = "  -2.23%  "
$ws.Range("D14").Value = "'672.58"
$ws.Range("E14").Value = "  +5.31%  "
$ws.Range("D15").Value = "3.880.60"
$ws.Range("E15").Value = "  +0.47%  "
$ws.Range("E16").Value = "  -0.56%  "
$ws.Range("D17").Value = "67.751.05"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("E18").Value = "  -0.59%  "
$ws.Range("D19").Value = "3.340.40"
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("E20").Value = "  -1.84%  "
$ws.Range("D21").Value = "'10.96"
$ws.Range("E21").Value = "  +0.43%  "
$ws.Range("E22").Value = "  -1.18%  "
$ws.Range("D23").Value = "'5.43"
$ws.Range("E23").Value = "  +8.55%  "
$ws.Range("D24").Value = "'17.10"
$ws.Range("E24").Value = "  -3.16%  "
$ws.Range("D25").Value = "'98.92"
$ws.Range("E25").Value = "  +1.53%  "
$ws.Range("E26").Value = "  -3.50%  "
$ws.Range("E27").Value = "  -4.03%  "
$ws.Range("E28").Value = "  -3.38%  "
$ws.Range("D29").Value = "'33.69"
$ws.Range("E29").Value = "  +2.77%  "
$ws.Range("D30").Value = "'8.43"
$ws.Range("E30").Value = "  -1.70%  "
$ws.Range("E31").Value = "  +10.56%  "
$ws.Range("D32").Value = "'572.79"
$ws.Range("E32").Value = "  -3.56%  "
$ws.Range("E33").Value = "  +0.25%  "
$ws.Range("E34").Value = "  +0.72%  "
$ws.Range("E35").Value = "  +0.16%  "
$ws.Range("D36").Value = "3.691.34"
$ws.Range("E36").Value = "  -6.24%  "
$ws.Range("D37").Value = "'56.63"
$ws.Range("E37").Value = "  +1.62%  "
$ws.Range("D38").Value = "'3.34"
$ws.Range("E38").Value = "  -5.39%  "
$ws.Range("D39").Value = "'34.44"
$ws.Range("E39").Value = "  +5.38%  "
$ws.Range("D40").Value = "'0.129"
$ws.Range("E40").Value = "  +0.55%  "
$ws.Range("E41").Value = "  -2.02%  "
$ws.Range("E42").Value = "  -4.62%  "
$ws.Range("D43").Value = "'3.32"
$ws.Range("E43").Value = "  -1.96%  "
$ws.Range("E44").Value = "  -1.22%  "
$ws.Range("E45").Value = "  -2.61%  "
$ws.Range("E46").Value = "  -1.92%  "
$ws.Range("E47").Value = "  +2.14%  "
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("E49").Value = "  -0.18%  "
$ws.Range("D50").Value = "'1.36"
$ws.Range("E50").Value = "  +1.18%  "
